$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows for Q4-Q7 (rows 6 through 9)
$ws.Rows("6:9").Delete()

# Update the remaining values for Q0 (row2), Q1 (row3), Q2 (row4), Q3 (row5)
$ws.Range("B2").Value = -0.001945245810387531
$ws.Range("C2").Value = 0.2935336163097054
$ws.Range("D2").Value = 0.1258855206509258
$ws.Range("E2").Value = 0.3548034958268109
$ws.Range("F2").Value = 0.3681914545215764

$ws.Range("B3").Value = -0.1174709636328633
$ws.Range("C3").Value = 0.1732624921278526
$ws.Range("D3").Value = 0.07702310541338803
$ws.Range("E3").Value = 0.2775303684525137
$ws.Range("F3").Value = 0.2650443948485162
$ws.Range("G3").Value = 10

$ws.Range("B4").Value = -0.09565720112800367
$ws.Range("C4").Value = 0.1894051810005109
$ws.Range("D4").Value = 0.09346882424652281
$ws.Range("E4").Value = 0.3057267149702865
$ws.Range("F4").Value = 0.3180915417653469
$ws.Range("G4").Value = 6

$ws.Range("B5").Value = -0.05251693463138896
$ws.Range("C5").Value = 0.05251693463138896
$ws.Range("D5").Value = 0.002890523596605528
$ws.Range("E5").Value = 0.05376358987833241
$ws.Range("F5").Value = 0.01627852410557828
$ws.Range("G5").Value = 2
